$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '20.077.31'
$ws.Range("E2").Value = '  -1.73%  '
$ws.Range("D3").Value = '1.421.01'
$ws.Range("E3").Value = '  -1.43%  '
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("E5").Value = '  -0.60%  '
$ws.Range("D6").Value = '''276.81'
$ws.Range("E6").Value = '  -0.46%  '
$ws.Range("D7").Value = '''0.3704'
$ws.Range("E7").Value = '  -1.37%  '
$ws.Range("D8").Value = '''0.3145'
$ws.Range("E8").Value = '  +2.56%  '
$ws.Range("D9").Value = '''39.67'
$ws.Range("E9").Value = '  -2.07%  '
$ws.Range("D10").Value = '''1.062'
$ws.Range("E10").Value = '  +4.97%  '
$ws.Range("E11").Value = '  -0.02%  '
$ws.Range("D12").Value = '''0.9988'
$ws.Range("E12").Value = '  -0.46%  '
$ws.Range("D13").Value = '''5.546'
$ws.Range("E13").Value = '  +3.25%  '
$ws.Range("D14").Value = '''17.96'
$ws.Range("E14").Value = '  +3.80%  '
$ws.Range("E15").Value = '  +1.01%  '
$ws.Range("D16").Value = '1.424.14'
$ws.Range("E16").Value = '  -1.30%  '
$ws.Range("E17").Value = '  +1.28%  '
$ws.Range("D18").Value = '''0.05717'
$ws.Range("E18").Value = '  -2.77%  '
$ws.Range("E19").Value = '  -0.54%  '
$ws.Range("D20").Value = '''71.76'
$ws.Range("E20").Value = '  -5.98%  '
$ws.Range("D21").Value = '''5.617'
$ws.Range("E21").Value = '  -1.97%  '
$ws.Range("D22").Value = '''14.89'
$ws.Range("E22").Value = '  +3.29%  '
$ws.Range("E23").Value = '  +1.81%  '
$ws.Range("E24").Value = '  -3.62%  '
$ws.Range("D25").Value = '20.134.08'
$ws.Range("D26").Value = '''2.296'
$ws.Range("E26").Value = '  +3.51%  '
$ws.Range("D27").Value = '''134.68'
$ws.Range("E27").Value = '  -5.85%  '
$ws.Range("E28").Value = '  +1.76%  '
$ws.Range("D29").Value = '1.582.18'
$ws.Range("E29").Value = '  -1.51%  '
$ws.Range("D30").Value = '''111.18'
$ws.Range("E30").Value = '  +1.28%  '
$ws.Range("D31").Value = '''3.963'
$ws.Range("E31").Value = '  +0.69%  '
$ws.Range("D32").Value = '''5.290'
$ws.Range("E32").Value = '  -2.74%  '
$ws.Range("D33").Value = '''0.8294'
$ws.Range("E33").Value = '  -8.46%  '
$ws.Range("D34").Value = '''0.07815'
$ws.Range("E34").Value = '  +0.83%  '
$ws.Range("D35").Value = '''1.481'
$ws.Range("E35").Value = '  +8.84%  '
$ws.Range("D36").Value = '''4.922'
$ws.Range("E36").Value = '  +4.18%  '
$ws.Range("D37").Value = '''0.05863'
$ws.Range("E37").Value = '  +3.74%  '
$ws.Range("D38").Value = '''7.933'
$ws.Range("E38").Value = '  -4.45%  '
$ws.Range("D39").Value = '''0.9961'
$ws.Range("E39").Value = '  -0.62%  '
$ws.Range("D40").Value = '''10.73'
$ws.Range("E40").Value = '  -0.85%  '
$ws.Range("D41").Value = '''0.02068'
$ws.Range("E41").Value = '  +1.02%  '
$ws.Range("D42").Value = '''1.110'
$ws.Range("E42").Value = '  -2.83%  '
$ws.Range("D43").Value = '''0.1878'
$ws.Range("E43").Value = '  -2.04%  '
$ws.Range("D44").Value = '''0.5359'
$ws.Range("E44").Value = '  +0.60%  '
$ws.Range("B45").Value = 'PancakeSwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D45").Value = '''3.553'
$ws.Range("E45").Value = '  -0.98%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = '''12.33'
$ws.Range("E46").Value = '  +1.91%  '
$ws.Range("D47").Value = '''118.16'
$ws.Range("E47").Value = '  +5.64%  '
$ws.Range("D48").Value = '''0.5249'
$ws.Range("E48").Value = '  +1.72%  '
$ws.Range("D49").Value = '''1.789'
$ws.Range("E49").Value = '  -0.21%  '
$ws.Range("E50").Value = '  -1.18%  '
$ws.Range("D51").Value = '''0.9971'
$ws.Range("E51").Value = '  -0.60%  '
